$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: values change only (row stays visible) ---
$ws.Range("M3").Value = 12.52
$ws.Range("N3").Value = 7.51
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 1
$ws.Range("U3").Value = 1

# --- Row 4: becomes hidden, values reset to 0 / stock adjustments ---
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 0
$ws.Range("U4").Value = 0
$ws.Rows(4).Hidden = $true

# --- Row 6: becomes hidden, values reset ---
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 0
$ws.Range("U6").Value = 0
$ws.Rows(6).Hidden = $true

# --- Row 8: becomes hidden, values reset ---
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("P8").Value = 2
$ws.Range("Q8").Value = 0
$ws.Range("U8").Value = 0
$ws.Rows(8).Hidden = $true

# --- Row 11: becomes hidden, values reset ---
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("P11").Value = 244
$ws.Range("Q11").Value = 0
$ws.Range("U11").Value = 0
$ws.Rows(11).Hidden = $true

# --- Row 12: already hidden, only L12 changes (unhide/rehide to avoid a
# row-height stamp being recorded while writing into a hidden row) ---
$ws.Rows(12).Hidden = $false
$ws.Range("L12").Value = -2
$ws.Rows(12).Hidden = $true

# --- Row 14: becomes hidden, values reset ---
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("P14").Value = 4
$ws.Range("Q14").Value = 0
$ws.Range("U14").Value = 0
$ws.Rows(14).Hidden = $true

# --- Row 16: becomes hidden, values reset ---
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("P16").Value = 6
$ws.Range("Q16").Value = 0
$ws.Range("U16").Value = 0
$ws.Rows(16).Hidden = $true

# --- Row 17: becomes hidden, values reset ---
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("P17").Value = 4
$ws.Range("Q17").Value = 0
$ws.Range("U17").Value = 0
$ws.Rows(17).Hidden = $true

# --- Row 20: becomes hidden, values reset ---
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("P20").Value = 4
$ws.Range("Q20").Value = 0
$ws.Range("U20").Value = 0
$ws.Rows(20).Hidden = $true

# --- Summary metrics ---
$ws.Range("C24").Value = 30
$ws.Range("C26").Value = "124.76€"
$ws.Range("C35").Value = -2
